$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Leading apostrophe forces Excel to store these as text (matching the
# original inline-string cell type) instead of auto-converting
# numeric-looking strings (e.g. '237.05', '0.05650') into numbers,
# which would silently change their representation (e.g. drop trailing
# zeros or merge '1.858.70' oddly).

# Row 2
$ws.Range("D2").Value = "'30.213.09"
$ws.Range("E2").Value = "'  +0.17%  "
# Row 3
$ws.Range("D3").Value = "'1.858.70"
$ws.Range("E3").Value = "'  -0.07%  "
# Row 4
$ws.Range("E4").Value = "'  +0.09%  "
# Row 5
$ws.Range("D5").Value = "'237.05"
$ws.Range("E5").Value = "'  +1.30%  "
# Row 6
$ws.Range("E6").Value = "'  +0.07%  "
# Row 7
$ws.Range("D7").Value = "'0.4677"
$ws.Range("E7").Value = "'  +0.15%  "
# Row 8
$ws.Range("D8").Value = "'0.2865"
$ws.Range("E8").Value = "'  +1.28%  "
# Row 9
$ws.Range("D9").Value = "'0.06534"
$ws.Range("E9").Value = "'  +0.12%  "
# Row 10
$ws.Range("D10").Value = "'21.92"
$ws.Range("E10").Value = "'  +5.40%  "
# Row 11
$ws.Range("D11").Value = "'0.07942"
$ws.Range("E11").Value = "'  +1.18%  "
# Row 12
$ws.Range("D12").Value = "'97.15"
$ws.Range("E12").Value = "'  +0.78%  "
# Row 13
$ws.Range("D13").Value = "'1.869.52"
$ws.Range("E13").Value = "'  +0.45%  "
# Row 14
$ws.Range("D14").Value = "'5.177"
$ws.Range("E14").Value = "'  +1.34%  "
# Row 15
$ws.Range("D15").Value = "'0.6803"
$ws.Range("E15").Value = "'  +1.60%  "
# Row 16
$ws.Range("D16").Value = "'267.42"
$ws.Range("E16").Value = "'  -3.77%  "
# Row 17
$ws.Range("D17").Value = "'30.203.78"
$ws.Range("E17").Value = "'  +0.09%  "
# Row 18
$ws.Range("D18").Value = "'13.66"
$ws.Range("E18").Value = "'  +8.17%  "
# Row 19
$ws.Range("E19").Value = "'  +0.22%  "
# Row 20
$ws.Range("D20").Value = "'0.000007372"
$ws.Range("E20").Value = "'  +1.84%  "
# Row 21
$ws.Range("D21").Value = "'2.112.94"
$ws.Range("E21").Value = "'  +0.08%  "
# Row 22
$ws.Range("D22").Value = "'5.321"
$ws.Range("E22").Value = "'  -2.58%  "
# Row 24
$ws.Range("D24").Value = "'6.198"
$ws.Range("E24").Value = "'  +0.93%  "
# Row 25
$ws.Range("D25").Value = "'167.24"
$ws.Range("E25").Value = "'  +1.54%  "
# Row 26
$ws.Range("D26").Value = "'9.219"
$ws.Range("E26").Value = "'  -0.86%  "
# Row 27
$ws.Range("D27").Value = "'18.88"
$ws.Range("E27").Value = "'  -0.45%  "
# Row 28
$ws.Range("D28").Value = "'1.951"
$ws.Range("E28").Value = "'  +2.22%  "
# Row 29
$ws.Range("D29").Value = "'1.386"
$ws.Range("E29").Value = "'  +1.74%  "
# Row 30
$ws.Range("D30").Value = "'0.09836"
$ws.Range("E30").Value = "'  +3.03%  "
# Row 31
$ws.Range("B31").Value = "'PancakeSwap"
$ws.Range("C31").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.485"
$ws.Range("E31").Value = "'  +1.28%  "
# Row 32
$ws.Range("B32").Value = "'Filecoin"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.367"
$ws.Range("E32").Value = "'  -0.83%  "
# Row 33
$ws.Range("D33").Value = "'4.056"
$ws.Range("E33").Value = "'  -0.99%  "
# Row 34
$ws.Range("D34").Value = "'0.04718"
$ws.Range("E34").Value = "'  +0.80%  "
# Row 35
$ws.Range("D35").Value = "'1.131"
$ws.Range("E35").Value = "'  +3.19%  "
# Row 36
$ws.Range("E36").Value = "'  +0.29%  "
# Row 37
$ws.Range("D37").Value = "'2.708"
$ws.Range("E37").Value = "'  -0.10%  "
# Row 38
$ws.Range("D38").Value = "'0.01877"
$ws.Range("E38").Value = "'  +0.29%  "
# Row 39
$ws.Range("D39").Value = "'2.629"
$ws.Range("E39").Value = "'  +4.33%  "
# Row 40
$ws.Range("D40").Value = "'75.23"
$ws.Range("E40").Value = "'  +3.47%  "
# Row 41
$ws.Range("D41").Value = "'6.243"
$ws.Range("E41").Value = "'  -1.22%  "
# Row 42
$ws.Range("D42").Value = "'1.946"
$ws.Range("E42").Value = "'  +0.58%  "
# Row 43
$ws.Range("D43").Value = "'0.8499"
$ws.Range("E43").Value = "'  +0.45%  "
# Row 44
$ws.Range("B44").Value = "'TheSandbox"
$ws.Range("C44").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4163"
$ws.Range("E44").Value = "'  -0.07%  "
# Row 45
$ws.Range("B45").Value = "'PaxDollar"
$ws.Range("C45").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9999"
$ws.Range("E45").Value = "'  +0.01%  "
# Row 46
$ws.Range("D46").Value = "'103.21"
$ws.Range("E46").Value = "'  -0.39%  "
# Row 47
$ws.Range("D47").Value = "'956.26"
$ws.Range("E47").Value = "'  -3.89%  "
# Row 48
$ws.Range("D48").Value = "'7.171"
$ws.Range("E48").Value = "'  +0.54%  "
# Row 49
$ws.Range("D49").Value = "'9.259"
$ws.Range("E49").Value = "'  -0.01%  "
# Row 50
$ws.Range("D50").Value = "'34.11"
$ws.Range("E50").Value = "'  +0.47%  "
# Row 51
$ws.Range("D51").Value = "'0.05650"
$ws.Range("E51").Value = "'  +0.75%  "
